# Applies the updated enrollment counts to the "Inscricoes" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 4
$ws.Range("E4").Value = 28
$ws.Range("F4").Value = 16
$ws.Range("H4").Value = 16

# Row 6
$ws.Range("E6").Value = 28

# Row 12
$ws.Range("E12").Value = 268

# Row 17
$ws.Range("E17").Value = 52

# Row 25
$ws.Range("E25").Value = 122
$ws.Range("F25").Value = 52
$ws.Range("H25").Value = 52

# Row 29
$ws.Range("E29").Value = 113
$ws.Range("F29").Value = 63
$ws.Range("H29").Value = 63

# Row 41
$ws.Range("F41").Value = 86
$ws.Range("H41").Value = 86

# Row 42
$ws.Range("E42").Value = 202
$ws.Range("F42").Value = 95
$ws.Range("H42").Value = 95

# Row 44
$ws.Range("E44").Value = 170

# Row 47
$ws.Range("E47").Value = 251

# Row 48
$ws.Range("E48").Value = 121
